$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Cells.Item(32, 6).Value = 76
$ws.Cells.Item(32, 7).Value = 563.16

# Row 36
$ws.Cells.Item(36, 6).Value = 104
$ws.Cells.Item(36, 7).Value = 2630.16

# Row 53
$ws.Cells.Item(53, 6).Value = 64
$ws.Cells.Item(53, 7).Value = 5089.28

# Row 55
$ws.Cells.Item(55, 2).Value = 65137.38

# Row 184
$ws.Cells.Item(184, 6).Value = 132
$ws.Cells.Item(184, 7).Value = 4843.08

# Row 186
$ws.Cells.Item(186, 2).Value = 80962.49

# Row 255
$ws.Cells.Item(255, 6).Value = 60
$ws.Cells.Item(255, 7).Value = 6858

# Row 256
$ws.Cells.Item(256, 6).Value = 1768
$ws.Cells.Item(256, 7).Value = 32708

# Row 261
$ws.Cells.Item(261, 2).Value = 43651.54

# Row 266
$ws.Cells.Item(266, 6).Value = 6
$ws.Cells.Item(266, 7).Value = 806.82

# Row 270
$ws.Cells.Item(270, 2).Value = 11212.31

# Row 283
$ws.Cells.Item(283, 6).Value = 96
$ws.Cells.Item(283, 7).Value = 5815.68

# Row 293
$ws.Cells.Item(293, 6).Value = 78
$ws.Cells.Item(293, 7).Value = 2445.3

# Row 306
$ws.Cells.Item(306, 6).Value = 56
$ws.Cells.Item(306, 7).Value = 3475.36

# Row 316
$ws.Cells.Item(316, 6).Value = 114
$ws.Cells.Item(316, 7).Value = 3655.98

# Row 326
$ws.Cells.Item(326, 6).Value = 65
$ws.Cells.Item(326, 7).Value = 3285.1

# Row 333
$ws.Cells.Item(333, 2).Value = 328206.92

# Row 335
$ws.Cells.Item(335, 6).Value = 181
$ws.Cells.Item(335, 7).Value = 31977.27

# Row 336
$ws.Cells.Item(336, 6).Value = 118
$ws.Cells.Item(336, 7).Value = 36089.12

# Row 353
$ws.Cells.Item(353, 6).Value = 24
$ws.Cells.Item(353, 7).Value = 7678.32

# Row 365
$ws.Cells.Item(365, 6).Value = 132
$ws.Cells.Item(365, 7).Value = 15078.36

# Row 367
$ws.Cells.Item(367, 6).Value = 96
$ws.Cells.Item(367, 7).Value = 10832.64

# Row 375
$ws.Cells.Item(375, 2).Value = 57802
$ws.Cells.Item(375, 6).Value = 81
$ws.Cells.Item(375, 7).Value = 11621.88

# Row 376
$ws.Cells.Item(376, 2).Value = 62791
$ws.Cells.Item(376, 6).Value = 43
$ws.Cells.Item(376, 7).Value = 6169.64

# Row 392
$ws.Cells.Item(392, 2).Value = 63040
$ws.Cells.Item(392, 6).Value = 69
$ws.Cells.Item(392, 7).Value = 7577.58

# Row 393
$ws.Cells.Item(393, 2).Value = 57870
$ws.Cells.Item(393, 6).Value = 0
$ws.Cells.Item(393, 7).Value = 0

# Row 400
$ws.Cells.Item(400, 6).Value = 10
$ws.Cells.Item(400, 7).Value = 2286.6

# Row 421
$ws.Cells.Item(421, 2).Value = 63043
$ws.Cells.Item(421, 3).Value = "HUL-Rexona Coconut&amp;Olive Oils 4X100G"
$ws.Cells.Item(421, 4).Value = 115.01
$ws.Cells.Item(421, 5).Value = 137.41
$ws.Cells.Item(421, 6).Value = 55
$ws.Cells.Item(421, 7).Value = 6325.55

# Row 422
$ws.Cells.Item(422, 2).Value = 53060
$ws.Cells.Item(422, 3).Value = "HUL-REXONA COCONUT&amp;OLIVE OILS 4x100g"
$ws.Cells.Item(422, 4).Value = 109.82
$ws.Cells.Item(422, 5).Value = 131.19
$ws.Cells.Item(422, 6).Value = 1
$ws.Cells.Item(422, 7).Value = 109.82

# Row 423
$ws.Cells.Item(423, 6).Value = 525
$ws.Cells.Item(423, 7).Value = 30817.5

# Row 426
$ws.Cells.Item(426, 6).Value = 98
$ws.Cells.Item(426, 7).Value = 11362.12

# Row 431
$ws.Cells.Item(431, 6).Value = 25
$ws.Cells.Item(431, 7).Value = 12359.5

# Row 444
$ws.Cells.Item(444, 6).Value = 19
$ws.Cells.Item(444, 7).Value = 7250.4

# Row 448
$ws.Cells.Item(448, 2).Value = 63007
$ws.Cells.Item(448, 6).Value = 1113
$ws.Cells.Item(448, 7).Value = 190690.29

# Row 449
$ws.Cells.Item(449, 2).Value = 57856
$ws.Cells.Item(449, 6).Value = 2
$ws.Cells.Item(449, 7).Value = 342.66

# Row 450
$ws.Cells.Item(450, 2).Value = 63008
$ws.Cells.Item(450, 6).Value = 615
$ws.Cells.Item(450, 7).Value = 92969.55

# Row 451
$ws.Cells.Item(451, 2).Value = 57857
$ws.Cells.Item(451, 6).Value = 3
$ws.Cells.Item(451, 7).Value = 453.51

# Row 461
$ws.Cells.Item(461, 6).Value = 199
$ws.Cells.Item(461, 7).Value = 11834.53

# Row 463
$ws.Cells.Item(463, 6).Value = 3
$ws.Cells.Item(463, 7).Value = 892.92

# Row 464
$ws.Cells.Item(464, 2).Value = 1348623.73

# Row 466
$ws.Cells.Item(466, 6).Value = 54
$ws.Cells.Item(466, 7).Value = 9913.32

# Row 481
$ws.Cells.Item(481, 2).Value = 87746.57

# Row 508
$ws.Cells.Item(508, 6).Value = 42
$ws.Cells.Item(508, 7).Value = 1551.9

# Row 511
$ws.Cells.Item(511, 6).Value = 245
$ws.Cells.Item(511, 7).Value = 34444.55

# Row 513
$ws.Cells.Item(513, 2).Value = 48145.6

# Row 561
$ws.Cells.Item(561, 6).Value = 189
$ws.Cells.Item(561, 7).Value = 5985.63

# Row 562
$ws.Cells.Item(562, 6).Value = 122
$ws.Cells.Item(562, 7).Value = 8780.34

# Row 566
$ws.Cells.Item(566, 6).Value = 254
$ws.Cells.Item(566, 7).Value = 15598.14

# Row 567
$ws.Cells.Item(567, 6).Value = 2
$ws.Cells.Item(567, 7).Value = 109.4

# Row 569
$ws.Cells.Item(569, 2).Value = 94708.34

# Row 576
$ws.Cells.Item(576, 6).Value = 7
$ws.Cells.Item(576, 7).Value = 17845.45

# Row 585
$ws.Cells.Item(585, 2).Value = 149017.45

# Row 611
$ws.Cells.Item(611, 6).Value = 124
$ws.Cells.Item(611, 7).Value = 6209.92

# Row 617
$ws.Cells.Item(617, 6).Value = 101
$ws.Cells.Item(617, 7).Value = 16091.32

# Row 618
$ws.Cells.Item(618, 6).Value = 95
$ws.Cells.Item(618, 7).Value = 2554.55

# Row 619
$ws.Cells.Item(619, 6).Value = 85
$ws.Cells.Item(619, 7).Value = 4114.85

# Row 620
$ws.Cells.Item(620, 6).Value = 99
$ws.Cells.Item(620, 7).Value = 15772.68

# Row 621
$ws.Cells.Item(621, 6).Value = 206
$ws.Cells.Item(621, 7).Value = 1985.84

# Row 627
$ws.Cells.Item(627, 2).Value = 101722.84

# Row 646
$ws.Cells.Item(646, 6).Value = 4
$ws.Cells.Item(646, 7).Value = 644.8

# Row 648
$ws.Cells.Item(648, 6).Value = 67
$ws.Cells.Item(648, 7).Value = 7388.09

# Row 649
$ws.Cells.Item(649, 2).Value = 20135.05

# Row 682
$ws.Cells.Item(682, 6).Value = 46
$ws.Cells.Item(682, 7).Value = 1335.38

# Row 696
$ws.Cells.Item(696, 2).Value = 48192.68

# Row 797
$ws.Cells.Item(797, 6).Value = 33
$ws.Cells.Item(797, 7).Value = 4021.05

# Row 798
$ws.Cells.Item(798, 6).Value = 121
$ws.Cells.Item(798, 7).Value = 12576.74

# Row 801
$ws.Cells.Item(801, 6).Value = 56
$ws.Cells.Item(801, 7).Value = 4641.84

# Row 803
$ws.Cells.Item(803, 6).Value = 74
$ws.Cells.Item(803, 7).Value = 6133.86

# Row 804
$ws.Cells.Item(804, 2).Value = 84482.52

# Row 871
$ws.Cells.Item(871, 6).Value = 233
$ws.Cells.Item(871, 7).Value = 8686.24

# Row 879
$ws.Cells.Item(879, 2).Value = 33116.06

# Row 954
$ws.Cells.Item(954, 6).Value = 1622
$ws.Cells.Item(954, 7).Value = 264564.42

# Row 956
$ws.Cells.Item(956, 6).Value = 264
$ws.Cells.Item(956, 7).Value = 38187.6

# Row 957
$ws.Cells.Item(957, 6).Value = 25
$ws.Cells.Item(957, 7).Value = 953.5

# Row 959
$ws.Cells.Item(959, 6).Value = 77
$ws.Cells.Item(959, 7).Value = 11385.22

# Row 962
$ws.Cells.Item(962, 2).Value = 359445.39

# Row 967
$ws.Cells.Item(967, 2).Value = 5714395.92

# Row 968
$ws.Cells.Item(968, 2).Value = 5714395.92
